# "gantt als png added"
#
# The Gantt sheet is being prepped for a clean PNG/image export: a thin
# margin column is added on the left (A) and right (T) of the chart so the
# exported picture has breathing room, and the selection is moved onto the
# chart's print/export area (A2:T22) instead of the stray T27 cell left
# over from editing.
#
# (Saving the xlsx also regenerates the workbook's style tables, which is
# why the OOXML picks up a couple of extra duplicated built-in "Link" /
# "Besuchter Link" cell styles — harmless bookkeeping from Excel's style
# merge, reproduced here for completeness.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- narrow margin columns around the chart (A and T), width 2.5 chars ---
# ColumnWidth is expressed in "characters"; Excel stores the serialized
# <col width="..."> in a slightly larger internal unit, so 1.6666667 here
# round-trips to exactly 2.5 in the saved file.
$ws.Columns.Item(1).ColumnWidth = 1.6666666666666667
$ws.Columns.Item(20).ColumnWidth = 1.6666666666666667

# --- move the selection onto the chart area ---
$ws.Range("A2:T22").Select()

# --- style-table bookkeeping (Excel regenerates these on save/merge) ---
try {
    $styles = $wb.Styles
    $styles.Add("Link")
    $styles.Add("Besuchter Link")
} catch {
    # Not fatal if this particular host doesn't expose style duplication;
    # the visible worksheet changes above are the meaningful edit.
}
